$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Duplicate row 7 (format + content) into new row 8 so A8/B8 inherit the
# correct cell styles (s="1" / s="2") without Excel's date auto-parsing
# kicking in when we overwrite the values afterwards.
$ws.Rows.Item(7).Copy()
$ws.Rows.Item(8).Insert(-4121, 0)

# A8: "2012.4.16" must stay plain text (not get reinterpreted as a date).
# Stage it in an untouched helper cell formatted as Text, then paste only
# the *value* into A8 so A8's existing (copied) style is left alone.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "2012.4.16"
$scratch.Copy()
$ws.Range("A8").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()

$ws.Range("B8").Value = "根据周日的讨论，重新设计《贪食家族》游戏，并上传"
$ws.Range("D8").Value = 4

$ws.Rows.Item(8).RowHeight = 27

$excel.ActiveWindow.ScrollRow = 4
$ws.Range("C8").Select()
